$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header) -------------------------------------------------
# Three brand-new columns (Area AdaptiveDecoder / Area Matrix /
# Effective Area) are inserted right after "Area (mm2)" (col D), so the
# old E:G (Power/Clock ns/Clock MHz) block is rebuilt one column further
# right (H:J) and the trailing Filename column moves from I to L.
$ws.Range("E1").Value = "Area AdaptiveDecoder"
$ws.Range("F1").Value = "Area Matrix"
$ws.Range("G1").Value = "Effective Area"
$ws.Range("H1").Value = "Power (mW)"
$ws.Range("I1").Value = "Clock (ns)"
$ws.Range("J1").Value = "Clock (MHz) "
$ws.Range("L1").Value = "Filename"

# --- Row 2 ------------------------------------------------------------
# No Area AdaptiveDecoder / Area Matrix / Effective Area data for this
# build, so those cells stay blank - clear whatever old content used to
# live in the old E2:G2 (Power/Clock) slots before it gets rewritten.
$ws.Range("E2:G2").ClearContents()
$ws.Range("H2").Value = 57.8
$ws.Range("I2").Value = 100
$ws.Range("J2").Formula = "=1000/I2"
$ws.Range("L2").Value = "build-dc-2014-05-05_19-23"

# --- Row 3 --------------------------------------------------------
$ws.Range("D3").Value = 0.47199999999999998
$ws.Range("E3").Value = 0.092
$ws.Range("F3").Value = 0.169
$ws.Range("G3").Formula = "=D3+1*(E3+F3)"
$ws.Range("H3").Value = 51.4
$ws.Range("I3").Value = 28.5
$ws.Range("J3").Formula = "=1000/I3"
$ws.Range("L3").Value = "build-dc-2014-05-05_22-51"

# --- Column widths for the new / resized columns -----------------------
$ws.Range("E1").ColumnWidth = 21.5
$ws.Range("F1").ColumnWidth = 12.333333333333334
$ws.Range("G1").ColumnWidth = 14.666666666666666
$ws.Range("H1").ColumnWidth = 11.166666666666666

# --- Selection moved to G3 -------------------------------------------
$ws.Range("G3").Select()
